$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '246.05'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '22.11'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.368'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05871'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '6.379'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9637'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1418'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.03547'
$ws.Range('E11').Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07352'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03033'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.458'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.09385'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001586'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.04805'
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0005900'
$ws.Range('E18').Value = '17OneONE'
$ws.Range('B19').Value = 'TigerCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.006285'
$ws.Range('E19').Value = '18TigerCashTCH'
$ws.Range('B20').Value = 'HotbitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.004087'
$ws.Range('E20').Value = '19HotbitTokenHTB'
$ws.Range('B21').Value = 'BitKan'
$ws.Range('C21').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0009864'
$ws.Range('E21').Value = '20BitKanKAN'
$ws.Range('B22').Value = 'NitroEx'
$ws.Range('C22').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.00009700'
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('B23').Value = 'LEO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.686'
$ws.Range('E23').Value = '22LEOLEO'
$ws.Range('B24').Value = 'BTSEToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.199'
$ws.Range('E24').Value = '23BTSETokenBTSE'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.3252'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0002471'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03852'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1076'
$ws.Range('E41').Value = '40BKEXTokenBKK'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.002440'
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.003038'
$ws.Range('E43').Value = '42KickTokenKICK'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.005745'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005660'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6511'
$ws.Range('E47').Value = '46CoinbaseStockTokenCOINBestin24h'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.03548'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002100'
